$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2050419, " 09/01/24", 45315, "Entregue", "SIM"),
  @(2050435, " 09/01/24", 45315, "Entregue", "SIM"),
  @(2050447, " 09/01/24", 45315, "Entregue", "SIM"),
  @(2050460, " 09/01/24", 45315, "Entregue", "SIM"),
  @(2050464, " 09/01/24", 45315, "Entregue", "SIM"),
  @(2057337, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057380, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057393, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057441, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057463, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057474, " 16/01/24", 45315, "Entregue", "SIM"),
  @(1243111, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2055535, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055536, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055555, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055556, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055561, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055562, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2057261, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057328, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057372, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057381, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057382, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2057383, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2059012, " 17/01/24", 45315, "Entregue", "SIM"),
  @(2055492, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055493, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2057284, " 16/01/24", 45315, "Entregue", "SIM"),
  @(2055511, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055532, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055559, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055560, " 15/01/24", 45315, "Entregue", "SIM"),
  @(2055567, " 15/01/24", 45315, "Entregue", "SIM")
)

$startRow = 1078
$r = $startRow
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[2]
  $ws.Cells.Item($r, 6).Value = $row[3]
  $ws.Cells.Item($r, 7).Value = $row[4]
  $r = $r + 1
}

$endRow = $r - 1
$dateRange = $ws.Range("C" + $startRow + ":E" + $endRow)
$dateRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output ("Last row written: " + $endRow)